# HoangMinhLe_Resume.docx - "Update after 1st month with MS" edit
#
# Changes:
#  1. Job title wording: "Software Development Engineers " -> "Software Engineer "
#     for all three occurrences (current MS role + two MS/Amazon internships).
#  2. Add " - 2 months" after "(October 2016 - Present)" for the current role,
#     and move the (hidden) _GoBack bookmark there.
#  3. Replace the three bullet points under the current Microsoft role.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Software Development Engineers " -> "Software Engineer " (3 occurrences)
# ---------------------------------------------------------------------------
$searchStart = 0
for ($i = 0; $i -lt 3; $i++) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("Software Development Engineers ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found) {
        break
    }
    $rng.Text = "Software Engineer "
    $searchStart = $rng.End
}

# ---------------------------------------------------------------------------
# 2. "(October 2016 - Present)" -> add " - 2 months" and move the _GoBack
#    bookmark to the end of that paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute("Present)", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found) {
    $insertPoint = $d.Range($rng.End, $rng.End)
    $insertPoint.InsertAfter(" " + [char]0x2013 + " 2 months")
    $insertPoint.Font.NameFarEast = "Times New Roman"
    $insertPoint.Font.NameBi = "Arial"
    $insertPoint.Font.Italic = 1
    $insertPoint.Font.ItalicBi = 1
    $insertPoint.Font.Color = 0

    # Remove the old _GoBack bookmark (currently after "Updating parsing logic.")
    $gb = $d.Bookmarks.Item("_GoBack")
    $gb.Delete()

    # Re-create it right after the newly inserted " - 2 months" text.
    $afterMonths = $d.Range($insertPoint.End, $insertPoint.End)
    $d.Bookmarks.Add("_GoBack", $afterMonths)
}

# ---------------------------------------------------------------------------
# 3. Update the three bullet points for the current Microsoft role.
# ---------------------------------------------------------------------------
$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute("Supporting migration to new architecture.", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found) {
    $rng.Text = "Working remotely with a Redmond team on a content delivery system for Windows 10. "
}

$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute("Fixing functional tests. ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found) {
    $rng.Text = "Supporting migration to new architecture."
}

$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute("Updating parsing logic.", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found) {
    $rng.Text = "Reduced local deployment complexity for functional testing."
}
